$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'331.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.10%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'44.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.39%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.555"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.30%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08279"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'2.039"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'3.46%"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'3.37%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.1125"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-3.64%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1900"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'2.42%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'10.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-12.90%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09989"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.58%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04668"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.16%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1057"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.67%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001273"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.84%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.04112"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-2.43%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005917"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.43%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'-0.23%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'4.435"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'2.61%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'3.61%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.3352"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-3.54%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1385"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-2.33%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D23").Value = "'0.001301"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'3.74%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004410"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.83%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'7.52%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003741"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.25%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02791"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'7.74%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05748"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'4.46%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007620"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.76%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1420"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.33%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007558"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.18%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.001975"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-2.15%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008316"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-0.55%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00007043"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-0.70%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.14%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0005803"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.15%"
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = 'BOLO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D48").Value = "'0.003584"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'1.83%"
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = 'CoinbaseStockToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D49").Value = "'0.002524"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'9.64%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.14%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.14%"
$ws.Range("E51").Style = "Normal"
